$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '46.000.26'
$ws.Range('E2').Value = '  -1.55%  '

# Row 3
$ws.Range('D3').Value = '2.376.41'
$ws.Range('E3').Value = '  +2.82%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '301.17'
$ws.Range('E5').Value = '  -0.54%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '99.13'
$ws.Range('E6').Value = '  -4.18%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').Value = '  -0.97%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.511'
$ws.Range('E9').Value = '  -4.04%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.50'
$ws.Range('E10').Value = '  -7.27%  '

# Row 11
$ws.Range('E11').Value = '  -1.74%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.18'
$ws.Range('E12').Value = '  -3.33%  '

# Row 13
$ws.Range('E13').Value = '  -0.71%  '

# Row 14
$ws.Range('D14').Value = '2.740.03'
$ws.Range('E14').Value = '  +2.76%  '

# Row 15
$ws.Range('D15').Value = '2.382.96'
$ws.Range('E15').Value = '  +3.03%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.818'
$ws.Range('E16').Value = '  -1.07%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.68'
$ws.Range('E17').Value = '  -3.26%  '

# Row 18
$ws.Range('D18').Value = '45.904.39'
$ws.Range('E18').Value = '  -1.73%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.77'
$ws.Range('E19').Value = '  -4.92%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0959'
$ws.Range('E20').Value = '  +0.86%  '

# Row 21
$ws.Range('E21').Value = '  -2.02%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.59'
$ws.Range('E22').Value = '  +0.66%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '244.63'
$ws.Range('E23').Value = '  -1.78%  '

# Row 24
$ws.Range('E24').Value = '  -5.00%  '

# Row 25
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.06%  '

# Row 26
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  -2.45%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '39.27'
$ws.Range('E27').Value = '  -10.59%  '

# Row 28
$ws.Range('E28').Value = '  -3.37%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').Value = '  -2.10%  '

# Row 30
$ws.Range('E30').Value = '  +20.76%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '21.07'
$ws.Range('E31').Value = '  +4.19%  '

# Row 32
$ws.Range('E32').Value = '  +6.74%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.55'
$ws.Range('E33').Value = '  -4.51%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '146.95'
$ws.Range('E34').Value = '  +0.15%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0776'
$ws.Range('E35').Value = '  -4.01%  '

# Row 36
$ws.Range('E36').Value = '  -0.39%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.94'
$ws.Range('E37').Value = '  +7.39%  '

# Row 38
$ws.Range('E38').Value = '  -3.55%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.87'
$ws.Range('E39').Value = '  -5.88%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.92'
$ws.Range('E40').Value = '  -4.71%  '

# Row 41
$ws.Range('E41').Value = '  -2.35%  '

# Row 42
$ws.Range('E42').Value = '  -7.16%  '

# Row 43
$ws.Range('D43').Value = '1.928.23'
$ws.Range('E43').Value = '  +3.98%  '

# Row 44
$ws.Range('E44').Value = '  -0.02%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.07'
$ws.Range('E45').Value = '  +2.67%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.78'
$ws.Range('E46').Value = '  -10.37%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.40'
$ws.Range('E47').Value = '  +4.43%  '

# Row 48
$ws.Range('E48').Value = '  -6.03%  '

# Row 49
$ws.Range('D49').Value = '2.611.03'
$ws.Range('E49').Value = '  +2.73%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '97.84'
$ws.Range('E50').Value = '  +0.05%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '68.55'
$ws.Range('E51').Value = '  -8.58%  '
